# Adding additional plot scripts for direct/indirect jobs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avg Demand Scenario")

# New "Notes and assumptions" column header
$ws.Range("D1").Value = "Notes and assumptions"

# State/location for WTIV rows
$ws.Range("C11").Value = "TX"
$ws.Range("C12").Value = "PA"
$ws.Range("C13").Value = "AL"

# COD pushed out to 2028 for WTIV 2 and WTIV 3
$ws.Range("B12").Value = 2028
$ws.Range("B13").Value = 2028

# Notes/assumptions for WTIV rows
$ws.Range("D11").Value = "Keppel Amfels (Begin construction right after Charybdis)"
$ws.Range("D12").Value = "Philly shipyard (Clear order books in 2025)"
$ws.Range("D13").Value = "VT Halter (Clear order books in 2025)"

# Make "Avg Demand Scenario" the active sheet/tab with D14 selected
$ws.Activate() | Out-Null
$ws.Range("D14").Select() | Out-Null
